$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Column width adjustments.
#    The stored (raw OOXML) "width" attribute is ColumnWidth + 5/6, so to
#    land on a target raw width we assign (target - 0.8333333333333334).
# ---------------------------------------------------------------------
$colsToWidth8 = @("B","C","G","I","J","K","L","O","P","Q","V","X","AA","AB","AC","AD","AG","AH")
foreach ($colLetter in $colsToWidth8) {
    $ws.Columns($colLetter).ColumnWidth = 7.166666666666667
}
$ws.Columns("T").ColumnWidth = 8.166666666666666

# ---------------------------------------------------------------------
# 2) Update data rows 2-5 (row 2 only changes its timestamp; rows 3-5 get
#    entirely new measurement data), then delete the old row 6.
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = [double]45118.50694444445

$row3 = @('45118.51388888889','3.872','3.077','0','8.933','6.961','2.671','13.472','4.945','2.731','3.318','3.75','4.024','1.013','3.36','4.394','3','0.091','0.415','44.184','9.550000000000001','3.103','5.932','3.765','0.631','7.023','2.537','3.058','3.861','4.402','0.647','12.217','2.154','3.549')
$row4 = @('45118.52083333334','23.696','17.908','0.674','51.97','42.632','18.372','68.123','28.884','13.361','19.367','20.954','22.2','5.996','18.776','26.549','15.727','0.226','0.91','278.803','52.464','17.332','35.162','18.937','2.628','34.364','15.205','13.923','16.494','22.317','0.435','61.702','10.135','21.461')
$row5 = @('45118.52777777778','12.23','9.27','0.31','26.97','22.02','9.41','41.32','14.97','7.06','10.03','10.89','11.55','3.12','9.75','13.76','8.23','0.12','0.52','141.48','27.45','9','18.3','9.92','1.43','19.92','7.89','7.38','8.73','11.66','0.33','37.66','5.32','11.11')

for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = [double]$row3[$i]
}
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = [double]$row4[$i]
}
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = [double]$row5[$i]
}

# Remove the now-obsolete 6th data row (was row 6 in the source file).
$ws.Rows("6:6").Delete()
